# Commit: "detect loop in a linked list"
#
# The author filled in the "Done" column (C) with the placeholder text
# "<->" for a handful of previously-blank divider/gap cells throughout
# the sheet, and marked the three LinkedList questions around the
# "detect loop in a linked list" topic (rows 139-141) as "Yes" in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows that did not exist before (completely blank rows between
#     sections). Setting a cell value on a previously non-existent row
#     creates the row. These need the same centered style ("s=4") that the
#     rest of column C uses, so explicitly (re)apply center alignment.
$newGapCells = @("C54", "C99", "C137", "C175")
foreach ($addr in $newGapCells) {
    $cell = $ws.Range($addr)
    $cell.Value = "<->"
    $cell.HorizontalAlignment = -4108  # xlCenter
}

# --- Existing rows whose column-C cell was blank (but already styled) -
#     fill them with the same "<->" placeholder used everywhere else in
#     the column.
$blankFillCells = @(
    "C55", "C100", "C138", "C176",
    "C212", "C213",
    "C236", "C237",
    "C273", "C274",
    "C294", "C295",
    "C334", "C335",
    "C354", "C355",
    "C400", "C401",
    "C408", "C409",
    "C470", "C471"
)
foreach ($addr in $blankFillCells) {
    $ws.Range($addr).Value = "<->"
}

# --- Mark the "reverse linked list" / "reverse in groups" / "detect loop
#     in a linked list" rows as done.
$ws.Range("C139").Value = "Yes"
$ws.Range("C140").Value = "Yes"
$ws.Range("C141").Value = "Yes"

# --- View state: the author was scrolled/zoomed in on the LinkedList
#     section (around row 140) when they made this edit.
$win = $excel.ActiveWindow
$win.Zoom = 55
$ws.Range("F140").Select() | Out-Null
